# Insert a new data row at row 118 (pushing the existing rows 118-169 down
# to 119-170, and extending the used range to A1:R170), then populate the
# new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(118).Insert()

$ws.Range("A118").Value = 8
$ws.Range("B118").Value = "Terminal La Palmera de La Serena"
$ws.Range("C118").Value = "Coquimbo"
$ws.Range("D118").Value = 44572
$ws.Range("E118").Value = 4
$ws.Range("F118").Value = 100112021
$ws.Range("G118").Value = "Ají"
$ws.Range("H118").Value = "Americana (o)"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 500
$ws.Range("K118").Value = 15000
$ws.Range("L118").Value = 16000
$ws.Range("M118").Value = 15500
$ws.Range("N118").Value = "$/caja 15 kilos"
$ws.Range("O118").Value = "Provincia de Limarí"
$ws.Range("P118").Value = 1033
$ws.Range("Q118").Value = 15
$ws.Range("R118").Value = "Hortaliza"
